# Weekly update: insert a new price record as row 65 (Fruta / hortaliza, semanal)
# This pushes the existing rows 65:167 down to 66:168 and grows the used range
# to A1:T168.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a brand-new blank row above the current row 65 (shifts 65:167 -> 66:168)
$ws.Rows.Item(65).Insert()

# Populate the new row 65 with the latest week's observation
$ws.Range("A65").Value = 4
$ws.Range("B65").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C65").Value = "Los Lagos"
$ws.Range("D65").Value = 45128
$ws.Range("E65").Value = 10
$ws.Range("F65").Value = "Fruta"
$ws.Range("G65").Value = 100104
$ws.Range("H65").Value = "Frutos de pepita"
$ws.Range("I65").Value = 100104003
$ws.Range("J65").Value = "Membrillo"
$ws.Range("K65").Value = "Champion"
$ws.Range("L65").Value = "Primera"
$ws.Range("M65").Value = 200
$ws.Range("N65").Value = 14000
$ws.Range("O65").Value = 14000
$ws.Range("P65").Value = 14000
$ws.Range("Q65").Value = "`$/caja 18 kilos empedrada"
$ws.Range("R65").Value = "Región de O'Higgins"
$ws.Range("S65").Value = 778
$ws.Range("T65").Value = 18

# Match the date format used by the rest of the "Fecha" column (D)
$ws.Range("D65").NumberFormat = $ws.Range("D66").NumberFormat
